# Commit: "Loan RBI, Variable Instalments"
#
# On the "Repayment schedule" sheet a new (blank) column was inserted
# immediately before column N. That shifts the old N/O/P columns ("Late",
# "heading", "Outstanding") one slot to the right (-> O/P/Q), the new
# column N is given a custom width, the sheet's selection moves to I18,
# and "Repayment schedule" becomes the active sheet/tab (so the former
# active sheet, "NewLoanInput", is no longer the selected tab).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = 10.2

$ws.Activate() | Out-Null
$ws.Range("I18").Select() | Out-Null
